$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.713.20'
$ws.Range("E2").Value = '  -0.13%  '
$ws.Range("D3").Value = '1.636.69'
$ws.Range("E3").Value = '  -0.66%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = '''217.16'
$ws.Range("E5").Value = '  +0.47%  '
$ws.Range("E6").Value = '  -1.06%  '
$ws.Range("E7").Value = '  +0.12%  '
$ws.Range("D8").Value = '''0.250'
$ws.Range("E8").Value = '  -1.09%  '
$ws.Range("D9").Value = '''0.0622'
$ws.Range("E9").Value = '  -0.80%  '
$ws.Range("D10").Value = '''19.03'
$ws.Range("E10").Value = '  -0.85%  '
$ws.Range("D11").Value = '''0.0845'
$ws.Range("E11").Value = '  +0.33%  '
$ws.Range("D12").Value = '1.862.17'
$ws.Range("E12").Value = '  -0.81%  '
$ws.Range("D13").Value = '1.639.38'
$ws.Range("E13").Value = '  -0.79%  '
$ws.Range("E14").Value = '  -1.16%  '
$ws.Range("D15").Value = '''0.524'
$ws.Range("E15").Value = '  -1.51%  '
$ws.Range("D16").Value = '''64.38'
$ws.Range("E16").Value = '  -1.71%  '
$ws.Range("D17").Value = '26.674.53'
$ws.Range("E18").Value = '  -2.59%  '
$ws.Range("E19").Value = '  +0.19%  '
$ws.Range("D20").Value = '''210.52'
$ws.Range("E20").Value = '  -3.92%  '
$ws.Range("E21").Value = '  -1.00%  '
$ws.Range("D22").Value = '''6.18'
$ws.Range("E22").Value = '  -1.51%  '
$ws.Range("D23").Value = '''2.31'
$ws.Range("E23").Value = '  -2.13%  '
$ws.Range("D24").Value = '''9.24'
$ws.Range("E24").Value = '  -2.93%  '
$ws.Range("D25").Value = '''145.82'
$ws.Range("E25").Value = '  -0.26%  '
$ws.Range("E26").Value = '  +0.11%  '
$ws.Range("E27").Value = '  -2.25%  '
$ws.Range("D28").Value = '''7.06'
$ws.Range("E28").Value = '  -0.73%  '
$ws.Range("D29").Value = '''15.50'
$ws.Range("E29").Value = '  -1.49%  '
$ws.Range("E30").Value = '  -2.72%  '
$ws.Range("D31").Value = '''1.19'
$ws.Range("E31").Value = '  +0.73%  '
$ws.Range("E32").Value = '  -0.39%  '
$ws.Range("D33").Value = '''2.97'
$ws.Range("E33").Value = '  -1.66%  '
$ws.Range("D34").Value = '1.272.45'
$ws.Range("E34").Value = '  -0.79%  '
$ws.Range("E35").Value = '  -1.61%  '
$ws.Range("E36").Value = '  +0.39%  '
$ws.Range("E37").Value = '  -2.17%  '
$ws.Range("D38").Value = '''0.525'
$ws.Range("E38").Value = '  -2.01%  '
$ws.Range("E39").Value = '  -2.57%  '
$ws.Range("E40").Value = '  +0.07%  '
$ws.Range("D41").Value = '''0.802'
$ws.Range("E41").Value = '  -1.58%  '
$ws.Range("E42").Value = '  -2.57%  '
$ws.Range("D43").Value = '1.772.79'
$ws.Range("E43").Value = '  -0.82%  '
$ws.Range("E44").Value = '  -3.71%  '
$ws.Range("D45").Value = '''91.26'
$ws.Range("E45").Value = '  -0.75%  '
$ws.Range("D46").Value = '''60.06'
$ws.Range("E46").Value = '  +0.62%  '
$ws.Range("E47").Value = '  -2.48%  '
$ws.Range("D48").Value = '0.0₆0101'
$ws.Range("E48").Value = '  -3.84%  '
$ws.Range("D49").Value = '''0.0518'
$ws.Range("E49").Value = '  +0.47%  '
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").Value = '''0.0961'
$ws.Range("E50").Value = '  -0.87%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Value = '''7.51'
$ws.Range("E51").Value = '  -3.49%  '
